$d = $word.ActiveDocument

# Locate the paragraph ending in "...file descriptors, redirect, and pipe"
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*file descriptors, redirect, and pipe*") {
        $targetIndex = $i
    }
}
if ($targetIndex -eq -1) {
    throw "Could not locate target paragraph"
}

# Create three placeholder empty paragraphs right after it
$r = $d.Paragraphs.Item($targetIndex).Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$emptyIndex = $targetIndex + 1
$r2 = $d.Paragraphs.Item($emptyIndex).Range
$r2.Collapse(0)
$r2.InsertParagraphAfter()

$headerIndex = $emptyIndex + 1
$r3 = $d.Paragraphs.Item($headerIndex).Range
$r3.Collapse(0)
$r3.InsertParagraphAfter()

$contentIndex = $headerIndex + 1

# 1) Blank spacer paragraph (justified)
$emptyXml = '<w:p><w:pPr><w:jc w:val="both"/></w:pPr></w:p>'
$d.Paragraphs.Item($emptyIndex).Range.InsertXML($emptyXml)

# 2) "October 22nd, 2022" header paragraph (bold, underline, "nd" superscript)
$hdrXml = '<w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:t>October 22</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/><w:vertAlign w:val="superscript"/></w:rPr><w:t>nd</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:t>, 2022</w:t></w:r></w:p>'
$d.Paragraphs.Item($headerIndex).Range.InsertXML($hdrXml)

# 3) New content paragraph describing the day's notes
$contentXml = '<w:p><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>Learned about &amp;&amp; operator and ‘cut’, ’sort’, ’</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>uniq</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>’, ‘</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>wc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>’, ‘grep’, and ‘apt-get’ programs. Some of these I learned in previous course, so they are only review.</w:t></w:r></w:p>'
$d.Paragraphs.Item($contentIndex).Range.InsertXML($contentXml)
